$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.272.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -4.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.938.23'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.58%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +2.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.932.47'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("E10").Value = '  -3.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.444'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.405.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.79'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +10.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.934.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '57.131.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '415.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.77%  '
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.680'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.28%  '
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.995'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.48'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("E33").Value = '  +9.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("E36").Value = '  -3.39%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '48.49'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0674'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.93%  '
$ws.Range("E39").Value = '  +6.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0346'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.107'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '375.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.633.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.239'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.22%  '
$ws.Range("E48").Value = '  +2.54%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("E51").Value = '  +0.48%  '
